# Insert a new column at the very left (becomes column A), shifting the
# existing A:E data to B:F, then populate the new column with row IDs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current column A. This shifts all existing
# columns (A-E) one place to the right (B-F) and shifts their formatting too.
$ws.Columns.Item(1).Insert()

# Header for the new ID column - match the style used by the other headers
# (bold / centered / bordered), which is already applied to row 1 cells
# that got shifted into B1:F1.
$ws.Cells.Item(1, 1).Value = "ID"
$ws.Cells.Item(1, 2).Copy() | Out-Null
$ws.Cells.Item(1, 1).PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Row labels for the new ID column (rows 2-25)
$ids = @{
    2  = "Hb 2"
    3  = "Hb 3"
    4  = "S 24"
    5  = "S 28"
    6  = "Hb 107"
    7  = "Hb 66"
    8  = "Hb 69"
    9  = "Hb 95"
    10 = "Hb 99"
    11 = "Hb 92"
    12 = "Hb 40"
    13 = "Hb 41"
    14 = "S 11"
    15 = "Hb 57"
    16 = "S 21"
    17 = "S 22"
    18 = "S 3"
    19 = "S 4"
    20 = "S 5"
    21 = "Hb 74"
    22 = "Hb 79"
    23 = "Hb 32"
    24 = "S 15"
    25 = "S 16"
}

foreach ($row in $ids.Keys) {
    $ws.Cells.Item($row, 1).Value = $ids[$row]
}
